$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right after the header row (row 1), pushing all
# existing data rows (old row 2 .. old row 63) down by two positions.
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).Insert()

# The inserted rows picked up the header row's formatting (bold/border); reset
# them to the plain, unstyled look used by the rest of the data rows, then
# reapply just the date format used in column D for data rows.
$ws.Rows.Item(2).ClearFormats()
$ws.Rows.Item(3).ClearFormats()
$ws.Cells.Item(2,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(3,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# New row 2: Angeleno / Primera / Región del Maule entry
$ws.Cells.Item(2,1).Value2 = 11
$ws.Cells.Item(2,2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item(2,3).Value2 = "Bíobío"
$ws.Cells.Item(2,4).Value2 = 44643
$ws.Cells.Item(2,5).Value2 = 8
$ws.Cells.Item(2,6).Value2 = "Fruta"
$ws.Cells.Item(2,7).Value2 = 100103
$ws.Cells.Item(2,8).Value2 = "Frutos de hueso (carozo)"
$ws.Cells.Item(2,9).Value2 = 100103002
$ws.Cells.Item(2,10).Value2 = "Ciruela"
$ws.Cells.Item(2,11).Value2 = "Angeleno"
$ws.Cells.Item(2,12).Value2 = "Primera"
$ws.Cells.Item(2,13).Value2 = 180
$ws.Cells.Item(2,14).Value2 = 8500
$ws.Cells.Item(2,15).Value2 = 9000
$ws.Cells.Item(2,16).Value2 = 8778
$ws.Cells.Item(2,17).Value2 = "$/bandeja 18 kilos granel"
$ws.Cells.Item(2,18).Value2 = "Región del Maule"
$ws.Cells.Item(2,19).Value2 = 488
$ws.Cells.Item(2,20).Value2 = 18

# New row 3: Angeleno / Segunda / Región del Maule entry
$ws.Cells.Item(3,1).Value2 = 11
$ws.Cells.Item(3,2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item(3,3).Value2 = "Bíobío"
$ws.Cells.Item(3,4).Value2 = 44643
$ws.Cells.Item(3,5).Value2 = 8
$ws.Cells.Item(3,6).Value2 = "Fruta"
$ws.Cells.Item(3,7).Value2 = 100103
$ws.Cells.Item(3,8).Value2 = "Frutos de hueso (carozo)"
$ws.Cells.Item(3,9).Value2 = 100103002
$ws.Cells.Item(3,10).Value2 = "Ciruela"
$ws.Cells.Item(3,11).Value2 = "Angeleno"
$ws.Cells.Item(3,12).Value2 = "Segunda"
$ws.Cells.Item(3,13).Value2 = 150
$ws.Cells.Item(3,14).Value2 = 6500
$ws.Cells.Item(3,15).Value2 = 7000
$ws.Cells.Item(3,16).Value2 = 6767
$ws.Cells.Item(3,17).Value2 = "$/bandeja 18 kilos granel"
$ws.Cells.Item(3,18).Value2 = "Región del Maule"
$ws.Cells.Item(3,19).Value2 = 376
$ws.Cells.Item(3,20).Value2 = 18
